$d = $word.ActiveDocument

# --- Step 1: locate the "Fast" paragraph (Solid State > Positives > Fast) and insert the new
#     block of paragraphs (Negatives/Expensive/.../Magnetic section) right after it, before the
#     two pre-existing blank paragraphs that separate it from the old trailing "Negatives" block. ---
$findRng = $d.Content
$found = $findRng.Find.Execute("Fast", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph 'Fast'"
}
$insertPos = $findRng.Paragraphs(1).Range.End
$insertAt = $d.Range($insertPos, $insertPos)

$xmlFrag = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Negatives</w:t></w:r></w:p><w:p><w:r><w:t>Expensive</w:t></w:r></w:p><w:p><w:r><w:t>Limited amount of read/writes</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>2026</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Optical Storage</w:t></w:r></w:p><w:p><w:r><w:t>Used by CDS and DVDS</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Uses lasers to read binary </w:t></w:r></w:p><w:p><w:r><w:t>If a part of the disc scatters slight, it is called a pit and gives a binary value of 0</w:t></w:r></w:p><w:p><w:r><w:t>If a part of the disc reflects light back into the laser, it is called a land and gives a binary of 1</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>CDS</w:t></w:r></w:p><w:p><w:r><w:t>Uses optical storage</w:t></w:r></w:p><w:p><w:r><w:t>Advantage: Portable, Cost effective</w:t></w:r></w:p><w:p><w:r><w:t>Disadvantage: Low storage space, slow</w:t></w:r><w:r><w:t>, easily scratched</w:t></w:r></w:p><w:p><w:r><w:t>Use cases: Audio</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>DVDS</w:t></w:r></w:p><w:p><w:r><w:t>Uses optical storage</w:t></w:r></w:p><w:p><w:r><w:t>Advantages: Larger storage space than CDS,</w:t></w:r><w:r><w:t xml:space="preserve"> more durable, portable</w:t></w:r></w:p><w:p><w:r><w:t>Disadvantages</w:t></w:r><w:r><w:t xml:space="preserve">: Eventually breaks down with natural use, still not durable or </w:t></w:r><w:r><w:br/><w:t>Use cases: Videos</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Blu ray</w:t></w:r></w:p><w:p><w:r><w:t>Advantage: Holds over 5x the capacity as DVDS, more durable than the others</w:t></w:r></w:p><w:p><w:r><w:t>Disadvantages: Expensive to produce, still slow</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Magnetic</w:t></w:r></w:p><w:p><w:r><w:t>Checks whether a portion of the region is polarized, as the magnetic poles will aligned and the signal between poles can be read.</w:t></w:r></w:p><w:p><w:r><w:t>Significantly more storage than optical</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$insertAt.InsertXML($xmlFrag)

# --- Step 2: the old trailing heading that used to read "Negatives" (followed by "Expensive" /
#     "Limited amount of read/writes") becomes "Hard disk drives (HDD)"; the two paragraphs below
#     it are removed and replaced by a single blank paragraph. Locate the LAST "Negatives"
#     heading in the document (the one immediately followed by "Expensive"). ---
$lastNegStart = -1
$lastNegEnd = -1
$searchRng = $d.Content
while ($true) {
    $found = $searchRng.Find.Execute("Negatives", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $lastNegStart = $searchRng.Start
    $lastNegEnd = $searchRng.End
    $searchRng.Collapse(0)
    $searchRng.End = $d.Content.End
}
if ($lastNegStart -eq -1) {
    throw "Could not find trailing 'Negatives' heading"
}

$negRng = $d.Range($lastNegStart, $lastNegEnd)
$negRng.Text = "Hard disk drives (HDD)"

# Remove the "Expensive" paragraph (text + its paragraph mark) immediately following.
$afterHeadingPos = $lastNegStart + ([string]"Hard disk drives (HDD)").Length
$expRng = $d.Range($afterHeadingPos, $afterHeadingPos)
$expRng.MoveEndUntil(-1, 1)
$expParaEnd = $expRng.Paragraphs(1).Range.End
$delRng = $d.Range($afterHeadingPos, $expParaEnd)
$delRng.Delete()

# Clear the text of the final "Limited amount of read/writes" paragraph, leaving it blank
# (its paragraph mark is the very last mark in the document body, so it cannot be deleted).
$lastParaRng = $d.Range($afterHeadingPos, $afterHeadingPos)
$lastParaEnd = $lastParaRng.Paragraphs(1).Range.End
$clearRng = $d.Range($afterHeadingPos, $lastParaEnd - 1)
if ($clearRng.Start -lt $clearRng.End) {
    $clearRng.Text = ""
}

Write-Host "done"
